$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row for a third cluster (interned first so it sorts before the
# updated formula strings in the shared string table)
$ws.Range("A3").Value = "IgGIV"

# Update existing trait formulas
$ws.Range("B1").Value = "first_trait = 0.5 * H4N4 + H5N4"
$ws.Range("B2").Value = "second_trait = (H4N4F1 + H5N4F1) / H3N4F1"
$ws.Range("B3").Value = "third_trait = H5N4F1S2 * H4N4F1S1"

# Move selection to reflect the edited state
$ws.Range("B4").Select()
